$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain text so values like "160.90" or "64.867.52"
# are not auto-coerced into numbers (which would drop formatting / trailing zeros).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '64.867.52'
$ws.Range("E2").Value = '  -1.06%  '
$ws.Range("D3").Value = '3.431.31'
$ws.Range("E3").Value = '  -1.56%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '574.57'
$ws.Range("E5").Value = '  -1.30%  '
$ws.Range("D6").Value = '159.17'
$ws.Range("E6").Value = '  -1.19%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").Value = '3.432.30'
$ws.Range("E8").Value = '  -1.65%  '
$ws.Range("B9").Value = 'XRP'
$ws.Range("C9").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D9").Value = '0.585'
$ws.Range("E9").Value = '  -3.90%  '
$ws.Range("E10").Value = '  -0.94%  '
$ws.Range("D11").Value = '0.123'
$ws.Range("E11").Value = '  -2.39%  '
$ws.Range("D12").Value = '0.447'
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("D13").Value = '4.024.53'
$ws.Range("E13").Value = '  -1.69%  '
$ws.Range("E15").Value = '  -3.40%  '
$ws.Range("D16").Value = '27.81'
$ws.Range("E16").Value = '  -3.11%  '
$ws.Range("D17").Value = '64.865.51'
$ws.Range("E17").Value = '  -1.05%  '
$ws.Range("D18").Value = '3.434.83'
$ws.Range("E18").Value = '  -1.55%  '
$ws.Range("D19").Value = '6.37'
$ws.Range("E19").Value = '  -1.32%  '
$ws.Range("D20").Value = '13.91'
$ws.Range("E20").Value = '  -2.84%  '
$ws.Range("D21").Value = '381.64'
$ws.Range("E21").Value = '  -1.70%  '
$ws.Range("D22").Value = '7.98'
$ws.Range("E22").Value = '  -3.20%  '
$ws.Range("E23").Value = '  -0.84%  '
$ws.Range("E24").Value = '  +0.21%  '
$ws.Range("D25").Value = '72.04'
$ws.Range("E25").Value = '  -1.72%  '
$ws.Range("E26").Value = '  -4.10%  '
$ws.Range("D27").Value = '9.94'
$ws.Range("E27").Value = '  -2.39%  '
$ws.Range("E28").Value = '  -0.64%  '
$ws.Range("E29").Value = '  +0.12%  '
$ws.Range("E30").Value = '  +2.25%  '
$ws.Range("D31").Value = '6.14'
$ws.Range("E31").Value = '  -2.84%  '
$ws.Range("E32").Value = '  -2.55%  '
$ws.Range("D33").Value = '23.29'
$ws.Range("E33").Value = '  -1.66%  '
$ws.Range("D34").Value = '7.06'
$ws.Range("E34").Value = '  -1.76%  '
$ws.Range("E35").Value = '  +0.43%  '
$ws.Range("D36").Value = '160.90'
$ws.Range("E36").Value = '  -0.95%  '
$ws.Range("E37").Value = '  -2.18%  '
$ws.Range("D38").Value = '2.903.25'
$ws.Range("E38").Value = '  -5.25%  '
$ws.Range("E39").Value = '  -3.04%  '
$ws.Range("D40").Value = '6.80'
$ws.Range("E40").Value = '  +4.01%  '
$ws.Range("D41").Value = '26.36'
$ws.Range("E41").Value = '  -3.12%  '
$ws.Range("D42").Value = '4.58'
$ws.Range("E42").Value = '  +0.25%  '
$ws.Range("D43").Value = '43.26'
$ws.Range("E43").Value = '  +0.24%  '
$ws.Range("E44").Value = '  -1.80%  '
$ws.Range("D45").Value = '0.773'
$ws.Range("E45").Value = '  -0.77%  '
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("E47").Value = '  +2.44%  '
$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").Value = '317.65'
$ws.Range("E48").Value = '  -0.13%  '
$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").Value = '1.08'
$ws.Range("E49").Value = '  -3.32%  '
$ws.Range("D50").Value = '6.52'
$ws.Range("E50").Value = '  -3.14%  '
$ws.Range("E51").Value = '  -3.33%  '
